$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '26.532.49'
$ws.Range('E2').Value2 = '  +0.66%  '
$ws.Range('D3').Value2 = '1.729.62'
$ws.Range('E3').Value2 = '  +0.67%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value2 = '0.9994'
$ws.Range('E4').Value2 = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '245.62'
$ws.Range('E5').Value2 = '  +2.80%  '
$ws.Range('E6').Value2 = '  -0.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value2 = '0.4808'
$ws.Range('E7').Value2 = '  +1.76%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.2674'
$ws.Range('E8').Value2 = '  +1.41%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '0.06222'
$ws.Range('E9').Value2 = '  +0.22%  '
$ws.Range('D10').Value2 = '1.727.61'
$ws.Range('E10').Value2 = '  +0.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value2 = '0.07167'
$ws.Range('E11').Value2 = '  +1.48%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value2 = '15.70'
$ws.Range('E12').Value2 = '  +2.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '0.6175'
$ws.Range('E13').Value2 = '  +4.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value2 = '4.539'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '77.20'
$ws.Range('E15').Value2 = '  +1.37%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value2 = '0.9997'
$ws.Range('E16').Value2 = '  -0.07%  '
$ws.Range('D17').Value2 = '26.533.55'
$ws.Range('E17').Value2 = '  +0.68%  '
$ws.Range('E18').Value2 = '  -0.08%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value2 = '0.000006953'
$ws.Range('E19').Value2 = '  +2.22%  '
$ws.Range('E20').Value2 = '  +0.90%  '
$ws.Range('D21').Value2 = '1.949.62'
$ws.Range('E22').Value2 = '  -0.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '8.933'
$ws.Range('E23').Value2 = '  +2.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value2 = '5.288'
$ws.Range('E24').Value2 = '  -0.65%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '136.64'
$ws.Range('E25').Value2 = '  +1.04%  '
$ws.Range('E26').Value2 = '  +0.70%  '
$ws.Range('E27').Value2 = '  +2.47%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value2 = '1.405'
$ws.Range('E28').Value2 = '  -0.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value2 = '3.992'
$ws.Range('E30').Value2 = '  -0.11%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '0.08019'
$ws.Range('E31').Value2 = '  +3.57%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value2 = '3.720'
$ws.Range('E32').Value2 = '  +0.85%  '
$ws.Range('E33').Value2 = '  +3.32%  '
$ws.Range('B34').Value2 = 'HuobiToken'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value2 = '2.617'
$ws.Range('E34').Value2 = '  +0.12%  '
$ws.Range('B35').Value2 = 'ImmutableX'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value2 = '0.6368'
$ws.Range('E35').Value2 = '  +2.93%  '
$ws.Range('B36').Value2 = 'ARBITRUM'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '0.9961'
$ws.Range('E36').Value2 = '  +1.72%  '
$ws.Range('B37').Value2 = 'TrustWalletToken'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '0.9226'
$ws.Range('E37').Value2 = '  -0.80%  '
$ws.Range('B38').Value2 = 'RenderToken'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value2 = '2.092'
$ws.Range('E38').Value2 = '  +9.80%  '
$ws.Range('B39').Value2 = 'MXToken'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value2 = '2.403'
$ws.Range('E39').Value2 = '  -0.39%  '
$ws.Range('B40').Value2 = 'Quant'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value2 = '104.73'
$ws.Range('E40').Value2 = '  -7.60%  '
$ws.Range('B41').Value2 = 'PaxDollar'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value2 = '1.003'
$ws.Range('E41').Value2 = '  +0.28%  '
$ws.Range('B42').Value2 = 'VeChain'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '0.01507'
$ws.Range('E42').Value2 = '  +1.97%  '
$ws.Range('B43').Value2 = 'FraxShare'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '5.596'
$ws.Range('E43').Value2 = '  +4.78%  '
$ws.Range('B44').Value2 = 'TheSandbox'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '0.3902'
$ws.Range('E44').Value2 = '  +2.49%  '
$ws.Range('B45').Value2 = 'Aptos'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '6.968'
$ws.Range('E45').Value2 = '  +10.66%  '
$ws.Range('B46').Value2 = 'Algorand'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '0.1184'
$ws.Range('E46').Value2 = '  +1.32%  '
$ws.Range('B47').Value2 = 'Cronos'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '0.05337'
$ws.Range('E47').Value2 = '  +0.93%  '
$ws.Range('B48').Value2 = 'Elrond'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '31.02'
$ws.Range('E48').Value2 = '  +1.97%  '
$ws.Range('B49').Value2 = 'EnergySwap'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value2 = '7.889'
$ws.Range('E49').Value2 = '  +2.64%  '
$ws.Range('B50').Value2 = 'NEARProtocol'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value2 = '1.267'
$ws.Range('E50').Value2 = '  +4.05%  '
$ws.Range('B51').Value2 = 'Decentraland'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '0.3429'
$ws.Range('E51').Value2 = '  +1.71%  '
